$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "model_4_3_1"
$ws.Range("B2").Value = 0.5698552127148301
$ws.Range("C2").Value = -0.3505366604854934
$ws.Range("D2").Value = 0.6501323513154735
$ws.Range("E2").Value = 0.06432008856170524
$ws.Range("F2").Value = 0.4760434925556183
$ws.Range("G2").Value = 1.815451502799988
$ws.Range("H2").Value = 0.3715589046478271
$ws.Range("I2").Value = 1.13597297668457

$ws.Range("A3").Value = "model_4_3_3"
$ws.Range("B3").Value = 0.6171988609290131
$ws.Range("C3").Value = -0.3441298721942034
$ws.Range("D3").Value = 0.722389260662942
$ws.Range("E3").Value = 0.09782038628046019
$ws.Range("F3").Value = 0.4236480593681335
$ws.Range("G3").Value = 1.806839227676392
$ws.Range("H3").Value = 0.2948221564292908
$ws.Range("I3").Value = 1.095301628112793

$ws.Range("A4").Value = "model_4_3_4"
$ws.Range("B4").Value = 0.6202536346294371
$ws.Range("C4").Value = -0.3380965063458155
$ws.Range("D4").Value = 0.7331921212445274
$ws.Range("E4").Value = 0.1058034741937606
$ws.Range("F4").Value = 0.4202672839164734
$ws.Range("G4").Value = 1.798728704452515
$ws.Range("H4").Value = 0.2833495438098907
$ws.Range("I4").Value = 1.085609555244446

$ws.Range("A5").Value = "model_4_3_2"
$ws.Range("B5").Value = 0.6244103050877862
$ws.Range("C5").Value = -0.2608440465238098
$ws.Range("D5").Value = 0.7217856476565176
$ws.Range("E5").Value = 0.1463923217375137
$ws.Range("F5").Value = 0.4156670868396759
$ws.Range("G5").Value = 1.694882750511169
$ws.Range("H5").Value = 0.2954632043838501
$ws.Range("I5").Value = 1.036332249641418

$ws.Range("A6").Value = "model_4_3_6"
$ws.Range("B6").Value = 0.6247724682407161
$ws.Range("C6").Value = -0.3377694537102744
$ws.Range("D6").Value = 0.7356239878005664
$ws.Range("E6").Value = 0.1069963921063649
$ws.Range("F6").Value = 0.4152662754058838
$ws.Range("G6").Value = 1.79828941822052
$ws.Range("H6").Value = 0.2807669043540955
$ws.Range("I6").Value = 1.084161281585693

$ws.Range("A7").Value = "model_4_3_5"
$ws.Range("B7").Value = 0.6255831657277684
$ws.Range("C7").Value = -0.3351565701766013
$ws.Range("D7").Value = 0.7443160226019847
$ws.Range("E7").Value = 0.1121068909424809
$ws.Range("F7").Value = 0.4143691062927246
$ws.Range("G7").Value = 1.794776916503906
$ws.Range("H7").Value = 0.2715359926223755
$ws.Range("I7").Value = 1.077956914901733

$ws.Range("A8").Value = "model_4_3_7"
$ws.Range("B8").Value = 0.6269986228016381
$ws.Range("C8").Value = -0.3346017149713765
$ws.Range("D8").Value = 0.7452240634368317
$ws.Range("E8").Value = 0.1128067553481981
$ws.Range("F8").Value = 0.4128026366233826
$ws.Range("G8").Value = 1.794031143188477
$ws.Range("H8").Value = 0.2705716490745544
$ws.Range("I8").Value = 1.077107310295105

$ws.Range("A9").Value = "model_4_3_8"
$ws.Range("B9").Value = 0.6297030882182301
$ws.Range("C9").Value = -0.3150918456784761
$ws.Range("D9").Value = 0.7357072265101904
$ws.Range("E9").Value = 0.1203239599592709
$ws.Range("F9").Value = 0.4098095595836639
$ws.Range("G9").Value = 1.767804980278015
$ws.Range("H9").Value = 0.2806785106658936
$ws.Range("I9").Value = 1.067980885505676

$ws.Range("A10").Value = "model_4_3_9"
$ws.Range("B10").Value = 0.6325929001836089
$ws.Range("C10").Value = -0.2969822147445305
$ws.Range("D10").Value = 0.7301469247255448
$ws.Range("E10").Value = 0.1286503809524016
$ws.Range("F10").Value = 0.4066114127635956
$ws.Range("G10").Value = 1.74346125125885
$ws.Range("H10").Value = 0.2865835428237915
$ws.Range("I10").Value = 1.05787205696106

$ws.Range("A11").Value = "model_4_3_10"
$ws.Range("B11").Value = 0.6336619665321834
$ws.Range("C11").Value = -0.2931018286744997
$ws.Range("D11").Value = 0.7323691380840114
$ws.Range("E11").Value = 0.1318394464650635
$ws.Range("F11").Value = 0.4054282605648041
$ws.Range("G11").Value = 1.738245129585266
$ws.Range("H11").Value = 0.2842235565185547
$ws.Range("I11").Value = 1.054000377655029

$ws.Range("A12").Value = "model_4_3_12"
$ws.Range("B12").Value = 0.6354901965809923
$ws.Range("C12").Value = -0.2842455028699264
$ws.Range("D12").Value = 0.7329558765506824
$ws.Range("E12").Value = 0.1372729485169701
$ws.Range("F12").Value = 0.4034049510955811
$ws.Range("G12").Value = 1.726339936256409
$ws.Range("H12").Value = 0.2836004197597504
$ws.Range("I12").Value = 1.047403812408447

$ws.Range("A13").Value = "model_4_3_11"
$ws.Range("B13").Value = 0.6355201661232528
$ws.Range("C13").Value = -0.2851640621457179
$ws.Range("D13").Value = 0.7345488844236396
$ws.Range("E13").Value = 0.1373902679303906
$ws.Range("F13").Value = 0.4033717811107635
$ws.Range("G13").Value = 1.727574825286865
$ws.Range("H13").Value = 0.2819086611270905
$ws.Range("I13").Value = 1.047261238098145

$ws.Range("A14").Value = "model_4_3_14"
$ws.Range("B14").Value = 0.6363308789568021
$ws.Range("C14").Value = -0.2810746682160996
$ws.Range("D14").Value = 0.7342831170973005
$ws.Range("E14").Value = 0.1396777322547571
$ws.Range("F14").Value = 0.402474582195282
$ws.Range("G14").Value = 1.722077608108521
$ws.Range("H14").Value = 0.2821908891201019
$ws.Range("I14").Value = 1.04448413848877

$ws.Range("A15").Value = "model_4_3_13"
$ws.Range("B15").Value = 0.6364875824996012
$ws.Range("C15").Value = -0.2801014300189011
$ws.Range("D15").Value = 0.7343017205393851
$ws.Range("E15").Value = 0.1402559077961177
$ws.Range("F15").Value = 0.4023011326789856
$ws.Range("G15").Value = 1.720769286155701
$ws.Range("H15").Value = 0.2821711301803589
$ws.Range("I15").Value = 1.043782234191895

$ws.Range("A16").Value = "model_4_3_15"
$ws.Range("B16").Value = 0.6378673422754128
$ws.Range("C16").Value = -0.2735224566033623
$ws.Range("D16").Value = 0.7346249793860167
$ws.Range("E16").Value = 0.1442452343432046
$ws.Range("F16").Value = 0.4007741510868073
$ws.Range("G16").Value = 1.711925506591797
$ws.Range("H16").Value = 0.2818278670310974
$ws.Range("I16").Value = 1.038938999176025

$ws.Range("A17").Value = "model_4_3_17"
$ws.Range("B17").Value = 0.6382092762631491
$ws.Range("C17").Value = -0.2614061484674439
$ws.Range("D17").Value = 0.7185947573870605
$ws.Range("E17").Value = 0.1447491633884115
$ws.Range("F17").Value = 0.4003957211971283
$ws.Range("G17").Value = 1.695638298988342
$ws.Range("H17").Value = 0.2988519072532654
$ws.Range("I17").Value = 1.038327217102051

$ws.Range("A18").Value = "model_4_3_24"
$ws.Range("B18").Value = 0.6387924901697406
$ws.Range("C18").Value = -0.2558437362741839
$ws.Range("D18").Value = 0.7142335865910153
$ws.Range("E18").Value = 0.1462146949463606
$ws.Range("F18").Value = 0.399750292301178
$ws.Range("G18").Value = 1.688161134719849
$ws.Range("H18").Value = 0.3034834861755371
$ws.Range("I18").Value = 1.036547899246216

$ws.Range("A19").Value = "model_4_3_19"
$ws.Range("B19").Value = 0.639055149509521
$ws.Range("C19").Value = -0.258021214710882
$ws.Range("D19").Value = 0.7196718688246571
$ws.Range("E19").Value = 0.147176501152759
$ws.Range("F19").Value = 0.3994596004486084
$ws.Range("G19").Value = 1.691088199615479
$ws.Range("H19").Value = 0.2977080345153809
$ws.Range("I19").Value = 1.035380244255066

$ws.Range("A20").Value = "model_4_3_21"
$ws.Range("B20").Value = 0.6391210709626526
$ws.Range("C20").Value = -0.2547204585556124
$ws.Range("D20").Value = 0.7152644163512909
$ws.Range("E20").Value = 0.1472971597539907
$ws.Range("F20").Value = 0.3993866443634033
$ws.Range("G20").Value = 1.686650991439819
$ws.Range("H20").Value = 0.3023887276649475
$ws.Range("I20").Value = 1.035233736038208

$ws.Range("A21").Value = "model_4_3_23"
$ws.Range("B21").Value = 0.6393750724115946
$ws.Range("C21").Value = -0.252677931587221
$ws.Range("D21").Value = 0.7141686340390767
$ws.Range("E21").Value = 0.1480427457460599
$ws.Range("F21").Value = 0.3991055190563202
$ws.Range("G21").Value = 1.683905482292175
$ws.Range("H21").Value = 0.3035524487495422
$ws.Range("I21").Value = 1.034328579902649

$ws.Range("A22").Value = "model_4_3_18"
$ws.Range("B22").Value = 0.639437271812294
$ws.Range("C22").Value = -0.2560234252797626
$ws.Range("D22").Value = 0.7197438788419143
$ws.Range("E22").Value = 0.1483779854474702
$ws.Range("F22").Value = 0.3990366756916046
$ws.Range("G22").Value = 1.688402652740479
$ws.Range("H22").Value = 0.297631561756134
$ws.Range("I22").Value = 1.033921480178833

$ws.Range("A23").Value = "model_4_3_20"
$ws.Range("B23").Value = 0.6394558959291943
$ws.Range("C23").Value = -0.2534020508239254
$ws.Range("D23").Value = 0.7159474844452689
$ws.Range("E23").Value = 0.1483500082761545
$ws.Range("F23").Value = 0.3990160524845123
$ws.Range("G23").Value = 1.684878826141357
$ws.Range("H23").Value = 0.301663339138031
$ws.Range("I23").Value = 1.033955574035645

$ws.Range("A24").Value = "model_4_3_22"
$ws.Range("B24").Value = 0.6396067525842409
$ws.Range("C24").Value = -0.2516758757544562
$ws.Range("D24").Value = 0.7144447231633606
$ws.Range("E24").Value = 0.1487451737736621
$ws.Range("F24").Value = 0.3988490998744965
$ws.Range("G24").Value = 1.682558417320251
$ws.Range("H24").Value = 0.3032592535018921
$ws.Range("I24").Value = 1.033475756645203

$ws.Range("A25").Value = "model_4_3_16"
$ws.Range("B25").Value = 0.639901065312843
$ws.Range("C25").Value = -0.2533361979105349
$ws.Range("D25").Value = 0.7203939846627243
$ws.Range("E25").Value = 0.1502195761469
$ws.Range("F25").Value = 0.3985233902931213
$ws.Range("G25").Value = 1.684790134429932
$ws.Range("H25").Value = 0.2969411611557007
$ws.Range("I25").Value = 1.031685829162598

$ws.Range("A26").Value = "model_4_3_0"
$ws.Range("B26").Value = 0.6669208972525469
$ws.Range("C26").Value = 0.2440808261337821
$ws.Range("D26").Value = 0.9749321129576011
$ws.Range("E26").Value = 0.5465763416594562
$ws.Range("F26").Value = 0.3686204254627228
$ws.Range("G26").Value = 1.016140222549438
$ws.Range("H26").Value = 0.02662205696105957
$ws.Range("I26").Value = 0.5504841804504395
